# Aggiornamento dati fino al 6/03: aggiunta di 3 nuove righe (245-247)
# in fondo alla tabella, replicando formato/stile dell'ultima riga esistente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dati da aggiungere: data (seriale), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(245, 44319, 5, 14, 141.9590346785642),
    @(246, 44320, 0, 14, 141.9590346785642),
    @(247, 44321, 1, 14, 141.9590346785642)
)

$lastRow = 244

foreach ($data in $newRows) {
    $r = $data[0]

    # Copia formattazione (stile, bordi, allineamento) dell'ultima riga esistente
    # nella nuova riga, cosi' la colonna A mantiene lo stile data (s="2").
    $ws.Range("A$lastRow`:D$lastRow").Copy($ws.Range("A$r`:D$r"))

    $ws.Cells.Item($r, 1).Value = $data[1]
    $ws.Cells.Item($r, 2).Value = $data[2]
    $ws.Cells.Item($r, 3).Value = $data[3]
    $ws.Cells.Item($r, 4).Value = $data[4]

    $lastRow = $r
}
